$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.196.71"
$ws.Range("E2").Value = "  -2.37%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.373.99"
$ws.Range("E3").Value = "  -0.46%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.71"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.43"
$ws.Range("E6").Value = "  +7.73%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.372.00"
$ws.Range("E8").Value = "  -0.58%  "

$ws.Range("E9").Value = "  +0.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.58"
$ws.Range("E10").Value = "  +2.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.123"
$ws.Range("E11").Value = "  +1.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.387"
$ws.Range("E12").Value = "  +2.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.946.72"
$ws.Range("E13").Value = "  -0.87%  "

$ws.Range("E14").Value = "  +1.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000174"
$ws.Range("E15").Value = "  +0.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.374.12"
$ws.Range("E16").Value = "  -0.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.12"
$ws.Range("E17").Value = "  +1.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.437.38"
$ws.Range("E18").Value = "  -2.05%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.02"
$ws.Range("E19").Value = "  +6.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.79"
$ws.Range("E20").Value = "  +1.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.35"
$ws.Range("E21").Value = "  -1.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "377.13"
$ws.Range("E22").Value = "  +0.24%  "

$ws.Range("E23").Value = "  +1.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.510.63"
$ws.Range("E24").Value = "  -0.56%  "

$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.72"
$ws.Range("E26").Value = "  -2.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000118"
$ws.Range("E27").Value = "  +9.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.68"
$ws.Range("E28").Value = "  +20.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.75"
$ws.Range("E29").Value = "  +10.86%  "

$ws.Range("E30").Value = "  +0.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.10"
$ws.Range("E31").Value = "  +3.16%  "

$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.402.89"
$ws.Range("E35").Value = "  -0.53%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.39"
$ws.Range("E36").Value = "  +2.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.57"
$ws.Range("E37").Value = "  +5.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.96"
$ws.Range("E38").Value = "  +3.12%  "

$ws.Range("E39").Value = "  +3.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "162.19"
$ws.Range("E40").Value = "  -1.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0784"
$ws.Range("E41").Value = "  +3.09%  "

$ws.Range("E42").Value = "  +0.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.22"
$ws.Range("E43").Value = "  +12.92%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.42"
$ws.Range("E44").Value = "  +3.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.55"
$ws.Range("E45").Value = "  +0.30%  "

$ws.Range("E46").Value = "  -2.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.61"
$ws.Range("E47").Value = "  +2.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.49"
$ws.Range("E48").Value = "  +3.04%  "

$ws.Range("E49").Value = "  +3.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.96"
$ws.Range("E50").Value = "  +13.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.902"
$ws.Range("E51").Value = "  +5.50%  "

# Row 32/33 swap: PancakeSwap <-> Kaspa reorder with updated values
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.157"
$ws.Range("E32").Value = "  +4.63%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.17"
$ws.Range("E33").Value = "  +0.32%  "
